$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 409 (shifts existing rows 409:496 down to 410:497).
$ws.Rows.Item(409).Insert()

# Populate the new row 409 with data (same record as the former row 409 /
# new row 410, but with an additional price point: Fecha, Volumen,
# Precio máximo, Precio promedio ponderado and Precio $/Kg differ).
$ws.Cells.Item(409, 1).Value = 10
$ws.Cells.Item(409, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(409, 3).Value = "La Araucanía"
$ws.Cells.Item(409, 4).Value = 45173
$ws.Cells.Item(409, 5).Value = 9
$ws.Cells.Item(409, 6).Value = "Fruta"
$ws.Cells.Item(409, 7).Value = 100102
$ws.Cells.Item(409, 8).Value = "Cítricos"
$ws.Cells.Item(409, 9).Value = 100102006
$ws.Cells.Item(409, 10).Value = "Pomelo"
$ws.Cells.Item(409, 11).Value = "Start Ruby"
$ws.Cells.Item(409, 12).Value = "Primera"
$ws.Cells.Item(409, 13).Value = 160
$ws.Cells.Item(409, 14).Value = 14000
$ws.Cells.Item(409, 15).Value = 15000
$ws.Cells.Item(409, 16).Value = 14625
$ws.Cells.Item(409, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(409, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(409, 19).Value = 975
$ws.Cells.Item(409, 20).Value = 15
